$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "jhnv"
$ws.Range("E7").Value = "vgnbcnbc"
$ws.Range("C8").Value = "vcbcvgb"
$ws.Range("B4").Value = "fbbfxfg"

$ws.Range("B4").Select()
